# Label BOM items better.
# Renames several "Package"/"Description" values in the BOM sheet to more
# consistent / clearer labels, narrows column B, and moves the active
# selection to D16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Package column (C): clearer part/footprint labels -------------------

# Resistors: "R-W4" -> "R-1/4W"
$ws.Range("C2").Value  = "R-1/4W"
$ws.Range("C3").Value  = "R-1/4W"
$ws.Range("C4").Value  = "R-1/4W"
$ws.Range("C5").Value  = "R-1/4W"
$ws.Range("C6").Value  = "R-1/4W"
$ws.Range("C7").Value  = "R-1/4W"
$ws.Range("C8").Value  = "R-1/4W"
$ws.Range("C9").Value  = "R-1/4W"
$ws.Range("C10").Value = "R-1/4W"

# Ceramic caps: "C-5mm" -> "C-P5mm"
$ws.Range("C12").Value = "C-P5mm"
$ws.Range("C13").Value = "C-P5mm"

# Film cap: "C-5mm 5x7.2mm" -> "C-P5mm 5x7.2mm"
$ws.Range("C14").Value = "C-P5mm 5x7.2mm"

# Electrolytic cap: "E2.5-6.3" -> "E-P2.5mm 6.3x11.5mm"
$ws.Range("C15").Value = "E-P2.5mm 6.3x11.5mm"

# --- Description column (E): reorder wording for consistency -------------

# "Capacitor Polarized THT" -> "Electrolytic Capacitor THT"
$ws.Range("E15").Value = "Electrolytic Capacitor THT"

# "Capacitor Film THT" -> "Film Capacitor THT"
$ws.Range("E14").Value = "Film Capacitor THT"

# "Capacitor Ceramic THT" -> "Ceramic Capacitor THT"
$ws.Range("E12").Value = "Ceramic Capacitor THT"
$ws.Range("E13").Value = "Ceramic Capacitor THT"

# --- Cosmetic sheet changes ------------------------------------------------

# Narrow the "Value" column (B)
$ws.Columns.Item(2).ColumnWidth = 12.14

# Move the active selection
$ws.Range("D16").Select() | Out-Null
